# BGU-14 RegNLic: add new TypeEditors rows (BreachOfLawRecordInfo .. FinancialOversightAuthorityInfo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TypeEditors")
$ws.Activate()

$names = @(
    "BreachOfLawRecordInfo",
    "EducationRecordInfo",
    "EmploymentRecordInfo",
    "FinancialGuaranteeInfo",
    "IncomeOriginInfo",
    "IndebtnessInfo",
    "IndebtnessInfoBase",
    "LiquidatedEntityOwnershipInfo",
    "LoanInfo",
    "PaymentDeadlineInfo",
    "PaymentModeInfo",
    "ProfessionLicenseInfo",
    "SharesAcquisitionInfo",
    "BankAccountInfo",
    "ProfessionLicensingBodyInfo",
    "PublicationInfo",
    "PublishingHouseInfo",
    "UniversityOrCollegeInfo",
    "FinancialOversightAuthorityInfo"
)

$startRow = 33

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $name = $names[$i]

    $ws.Cells.Item($r, 1).Value = $name

    $ws.Cells.Item($r, 2).Formula = '="I"& TRIM(A' + $r + ') & "EditFormFactory"'
    $ws.Cells.Item($r, 3).Formula = '="public interface I"& TRIM(A' + $r + ') & "EditFormFactory : ITypeEditorFormFactoryBase { }"'
    $ws.Cells.Item($r, 4).Formula = '=A' + $r + '& "_Editor"'
    $ws.Cells.Item($r, 5).Formula = '="public class " & D' + $r + ' & " : GenericTypeEditor<"&A' + $r + '&"> { private " & B' + $r + ' & " _fact; protected override ITypeEditorFormFactoryBase TypeEditorFormFactory { get { if (_fact == null) _fact = TypeEditorsDispatcher.Container.Resolve<" &B' + $r + ' & ">(); return _fact; } }  }"'
    $ws.Cells.Item($r, 6).Formula = '=A' + $r + '& "EditFormFactoryBasic"'
    $ws.Cells.Item($r, 7).Formula = '="public class " &F' + $r + '& " : " & B' + $r + ' & " { public System.Windows.Forms.Form SpawnInstance() { return new DummyForm<" &A' + $r + '& " >(); } }"'
    $ws.Cells.Item($r, 8).Formula = '="cont.RegisterInstance<" & B' + $r + ' & ">(new " & F' + $r + ' & "(), new ContainerControlledLifetimeManager());"'
    $ws.Cells.Item($r, 9).Formula = '="[System.ComponentModel.Editor(typeof(BGU.DRPL.SignificantOwnership.Core.TypeEditors." &D' + $r + ' & "), typeof(System.Drawing.Design.UITypeEditor))]"'
}

$endRow = $startRow + $names.Count - 1

$ws.Range("G" + $endRow).Select()
